$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to text before writing, so decimal-looking
# strings like "326.36" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.182.79"
$ws.Range("E2").Value = "  -1.48%  "

$ws.Range("D3").Value = "1.838.50"
$ws.Range("E3").Value = "  -0.66%  "

$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").Value = "326.36"
$ws.Range("E5").Value = "  -2.88%  "

$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").Value = "0.4636"
$ws.Range("E7").Value = "  -0.75%  "

$ws.Range("D8").Value = "0.3871"
$ws.Range("E8").Value = "  -1.05%  "

$ws.Range("D9").Value = "0.07858"
$ws.Range("E9").Value = "  -0.60%  "

$ws.Range("D10").Value = "0.9633"
$ws.Range("E10").Value = "  -1.70%  "

$ws.Range("D11").Value = "22.05"
$ws.Range("E11").Value = "  -1.09%  "

$ws.Range("D12").Value = "1.823.29"
$ws.Range("E12").Value = "  -3.27%  "

$ws.Range("D13").Value = "5.691"
$ws.Range("E13").Value = "  -2.30%  "

$ws.Range("D14").Value = "6.896"
$ws.Range("E14").Value = "  -1.18%  "

$ws.Range("D15").Value = "0.06870"
$ws.Range("E15").Value = "  -1.13%  "

$ws.Range("D16").Value = "88.51"
$ws.Range("E16").Value = "  +0.88%  "

$ws.Range("E17").Value = "  +0.29%  "

$ws.Range("D18").Value = "0.000009942"
$ws.Range("E18").Value = "  -0.85%  "

$ws.Range("D19").Value = "16.73"
$ws.Range("E19").Value = "  -2.07%  "

$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.18%  "

$ws.Range("D21").Value = "28.182.35"
$ws.Range("E21").Value = "  -1.55%  "

$ws.Range("D22").Value = "5.304"
$ws.Range("E22").Value = "  -1.73%  "

$ws.Range("D23").Value = "11.04"
$ws.Range("E23").Value = "  -2.05%  "

$ws.Range("D24").Value = "2.094"
$ws.Range("E24").Value = "  -2.62%  "

$ws.Range("D25").Value = "2.079.66"
$ws.Range("E25").Value = "  -0.76%  "

$ws.Range("D26").Value = "154.42"
$ws.Range("E26").Value = "  +1.12%  "

$ws.Range("D27").Value = "19.18"
$ws.Range("E27").Value = "  -1.31%  "

$ws.Range("D28").Value = "5.722"
$ws.Range("E28").Value = "  -5.91%  "

$ws.Range("E29").Value = "  -2.42%  "

$ws.Range("D30").Value = "119.21"
$ws.Range("E30").Value = "  +1.42%  "

$ws.Range("D31").Value = "0.9370"
$ws.Range("E31").Value = "  -3.66%  "

$ws.Range("D32").Value = "0.09269"
$ws.Range("E32").Value = "  -0.91%  "

$ws.Range("D33").Value = "5.284"
$ws.Range("E33").Value = "  -1.55%  "

$ws.Range("D34").Value = "1.323"
$ws.Range("E34").Value = "  -1.95%  "

$ws.Range("D35").Value = "3.323"
$ws.Range("E35").Value = "  -4.28%  "

$ws.Range("D36").Value = "0.05836"
$ws.Range("E36").Value = "  -5.00%  "

$ws.Range("D37").Value = "0.02126"
$ws.Range("E37").Value = "  -3.45%  "

$ws.Range("D38").Value = "1.140"
$ws.Range("E38").Value = "  -2.73%  "

$ws.Range("D39").Value = "7.766"
$ws.Range("E39").Value = "  +0.88%  "

$ws.Range("D40").Value = "0.5599"
$ws.Range("E40").Value = "  -2.01%  "

$ws.Range("D41").Value = "9.918"
$ws.Range("E41").Value = "  -2.12%  "

$ws.Range("D42").Value = "0.1761"
$ws.Range("E42").Value = "  -1.86%  "

$ws.Range("D43").Value = "0.07333"
$ws.Range("E43").Value = "  +3.28%  "

$ws.Range("D44").Value = "11.66"
$ws.Range("E44").Value = "  -0.55%  "

$ws.Range("D45").Value = "0.5276"
$ws.Range("E45").Value = "  -1.82%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "2.141"
$ws.Range("E46").Value = "  -12.43%  "

$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "1.141"
$ws.Range("E47").Value = "  -8.20%  "

$ws.Range("D48").Value = "1.842"
$ws.Range("E48").Value = "  -3.52%  "

$ws.Range("D49").Value = "114.03"
$ws.Range("E49").Value = "  +0.87%  "

$ws.Range("D50").Value = "1.001"
$ws.Range("E50").Value = "  +0.19%  "

$ws.Range("D51").Value = "2.324"
$ws.Range("E51").Value = "  -1.11%  "

# Drop the temporary text number-format again so the cells keep their
# original (default) style, matching the source formatting.
$ws.Range("D2:D51").ClearFormats()